# Matriz de trazabilidad_Ordoñez_Rivas.xlsx - update
# - "Fecha de estado" (column I) values bumped forward for every requirement row (3-39)
# - Row 24 (RF-24 "Gestionar Vale") status (column H) flips from "Falta" to "Completado"
# - Row 11 height shrinks from 78.75 to 63
# - Selection/cursor left on Q37 instead of Q5:R6, and the frozen top-left scroll cell is cleared

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Column I: "Fecha de estado" -------------------------------------------------
# Old value 43994 (2020-06-12) -> 44122 (2020-10-18) for most rows
$rowsShiftTo44122 = @(3,4,5,6,7,8,9,10,11,12,13,14,15,20,21,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39)
foreach ($r in $rowsShiftTo44122) {
    $ws.Cells.Item($r, 9).Value = 44122
}

# Old value 43987 (2020-06-05) -> 44109 (2020-10-05)
$rowsShiftTo44109 = @(16,17,19)
foreach ($r in $rowsShiftTo44109) {
    $ws.Cells.Item($r, 9).Value = 44109
}

# Old value 43990 (2020-06-08) -> 44112 (2020-10-08)
$ws.Cells.Item(18, 9).Value = 44112

# Old value 43995 (2020-06-13) -> 44117 (2020-10-13)
$rowsShiftTo44117 = @(22,23)
foreach ($r in $rowsShiftTo44117) {
    $ws.Cells.Item($r, 9).Value = 44117
}

# --- Row 24 (RF-24): Estado "Falta" -> "Completado" ------------------------------
$ws.Range("H24").Value = "Completado"

# --- Row 11 height: 78.75 -> 63 ---------------------------------------------------
$ws.Rows.Item(11).RowHeight = 63

# --- Selection: leave the cursor on Q37 instead of Q5:R6 -------------------------
$ws.Range("Q37").Select()
